# Insert a new data row at row 67 (pushing existing rows 67..172 down to 68..173)
# and populate it with the new day's price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(67).Insert()

$ws.Cells.Item(67, 1).Value = 5
$ws.Cells.Item(67, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(67, 3).Value = "Maule"
$ws.Cells.Item(67, 4).Value = 45070
$ws.Cells.Item(67, 5).Value = 7
$ws.Cells.Item(67, 6).Value = "Fruta"
$ws.Cells.Item(67, 7).Value = 100108
$ws.Cells.Item(67, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(67, 9).Value = 100108002
$ws.Cells.Item(67, 10).Value = "Mango"
$ws.Cells.Item(67, 11).Value = "Sin especificar"
$ws.Cells.Item(67, 12).Value = "Primera"
$ws.Cells.Item(67, 13).Value = 248
$ws.Cells.Item(67, 14).Value = 7000
$ws.Cells.Item(67, 15).Value = 7000
$ws.Cells.Item(67, 16).Value = 7000
$ws.Cells.Item(67, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(67, 18).Value = "Perú"
$ws.Cells.Item(67, 19).Value = 1750
$ws.Cells.Item(67, 20).Value = 4
